$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 185
$ws.Range("A185").Value = ""
$ws.Range("B185").Value = "Robust μ-Synthesis Grid-Side Control for Inverter-Based Resources in Weak Grids"
$ws.Range("C185").Value = 2026
$ws.Range("D185").Value = "MDPI AG"
$ws.Range("E185").Value = "Energies"
$ws.Range("F185").Value = "Kim, Woo-Jung; Lee, Yu-Seok; Chun, Yeong-Han"
$ws.Range("G185").Value = ""
$ws.Range("H185").Value = "10.3390/en19040946"
$ws.Range("I185").Value = "https://doi.org/10.3390/en19040946"
$ws.Range("J185").Value = "Journal"
$ws.Range("K185").Value = "Inverter"
$ws.Range("L185").Value = "Experiment"
$ws.Range("M185").Value = "Contacts"
$ws.Range("N185").Value = ""
$ws.Range("O185").Value = ""
$ws.Range("P185").Value = ""
$ws.Range("Q185").Value = "Robust μ-Synthesis Grid-Side Control for Inverter-Based Resources in Weak Grids"
$ws.Range("R185").Value = "High"
$ws.Range("S185").Value = "'2026-02-12"
$ws.Range("T185").Value = ""
$ws.Range("A185").Style = "Normal"
$ws.Range("G185").Style = "Normal"
$ws.Range("N185").Style = "Normal"
$ws.Range("O185").Style = "Normal"
$ws.Range("P185").Style = "Normal"
$ws.Range("S185").Style = "Normal"
$ws.Range("T185").Style = "Normal"

# Row 186
$ws.Range("A186").Value = ""
$ws.Range("B186").Value = "Inductance and capacitance parasitic prediction thanks to data analysis applied to Si and SiC MOSFET wide frequency band characterization"
$ws.Range("C186").Value = 2026
$ws.Range("D186").Value = "Elsevier BV"
$ws.Range("E186").Value = "Microelectronics Reliability"
$ws.Range("F186").Value = "Vidal, P.-E.; Viné, G.; Baffreau, S.; Gopishetti, A.; Le, T.L."
$ws.Range("G186").Value = ""
$ws.Range("H186").Value = "10.1016/j.microrel.2026.116047"
$ws.Range("I186").Value = "https://doi.org/10.1016/j.microrel.2026.116047"
$ws.Range("J186").Value = "Journal"
$ws.Range("K186").Value = "n-FET"
$ws.Range("L186").Value = "Experiment"
$ws.Range("M186").Value = "Gate Stack"
$ws.Range("N186").Value = ""
$ws.Range("O186").Value = ""
$ws.Range("P186").Value = ""
$ws.Range("Q186").Value = "Inductance and capacitance parasitic prediction thanks to data analysis applied to Si and SiC MOSFET wide frequency band characterization"
$ws.Range("R186").Value = "High"
$ws.Range("S186").Value = "'2026-02-12"
$ws.Range("T186").Value = ""
$ws.Range("A186").Style = "Normal"
$ws.Range("G186").Style = "Normal"
$ws.Range("N186").Style = "Normal"
$ws.Range("O186").Style = "Normal"
$ws.Range("P186").Style = "Normal"
$ws.Range("S186").Style = "Normal"
$ws.Range("T186").Style = "Normal"

# Row 187
$ws.Range("A187").Value = ""
$ws.Range("B187").Value = "Inductance and capacitance parasitic prediction thanks to data analysis applied to Si and SiC MOSFET wide frequency band characterization"
$ws.Range("C187").Value = 2026
$ws.Range("D187").Value = "Elsevier BV"
$ws.Range("E187").Value = "Microelectronics Reliability"
$ws.Range("F187").Value = "Vidal, P.-E.; Viné, G.; Baffreau, S.; Gopishetti, A.; Le, T.L."
$ws.Range("G187").Value = ""
$ws.Range("H187").Value = "10.1016/j.microrel.2026.116047"
$ws.Range("I187").Value = "https://doi.org/10.1016/j.microrel.2026.116047"
$ws.Range("J187").Value = "Journal"
$ws.Range("K187").Value = "n-FET"
$ws.Range("L187").Value = "Experiment"
$ws.Range("M187").Value = "Gate Stack"
$ws.Range("N187").Value = ""
$ws.Range("O187").Value = ""
$ws.Range("P187").Value = ""
$ws.Range("Q187").Value = "Inductance and capacitance parasitic prediction thanks to data analysis applied to Si and SiC MOSFET wide frequency band characterization"
$ws.Range("R187").Value = "High"
$ws.Range("S187").Value = "'2026-02-12"
$ws.Range("T187").Value = ""
$ws.Range("A187").Style = "Normal"
$ws.Range("G187").Style = "Normal"
$ws.Range("N187").Style = "Normal"
$ws.Range("O187").Style = "Normal"
$ws.Range("P187").Style = "Normal"
$ws.Range("S187").Style = "Normal"
$ws.Range("T187").Style = "Normal"

# Row 188
$ws.Range("A188").Value = ""
$ws.Range("B188").Value = "Robust μ-Synthesis Grid-Side Control for Inverter-Based Resources in Weak Grids"
$ws.Range("C188").Value = 2026
$ws.Range("D188").Value = "MDPI AG"
$ws.Range("E188").Value = "Energies"
$ws.Range("F188").Value = "Kim, Woo-Jung; Lee, Yu-Seok; Chun, Yeong-Han"
$ws.Range("G188").Value = ""
$ws.Range("H188").Value = "10.3390/en19040946"
$ws.Range("I188").Value = "https://doi.org/10.3390/en19040946"
$ws.Range("J188").Value = "Journal"
$ws.Range("K188").Value = "Inverter"
$ws.Range("L188").Value = "Experiment"
$ws.Range("M188").Value = "Contacts"
$ws.Range("N188").Value = ""
$ws.Range("O188").Value = ""
$ws.Range("P188").Value = ""
$ws.Range("Q188").Value = "Robust μ-Synthesis Grid-Side Control for Inverter-Based Resources in Weak Grids"
$ws.Range("R188").Value = "High"
$ws.Range("S188").Value = "'2026-02-12"
$ws.Range("T188").Value = ""
$ws.Range("A188").Style = "Normal"
$ws.Range("G188").Style = "Normal"
$ws.Range("N188").Style = "Normal"
$ws.Range("O188").Style = "Normal"
$ws.Range("P188").Style = "Normal"
$ws.Range("S188").Style = "Normal"
$ws.Range("T188").Style = "Normal"

# Row 189
$ws.Range("A189").Value = ""
$ws.Range("B189").Value = "Robust μ-Synthesis Grid-Side Control for Inverter-Based Resources in Weak Grids"
$ws.Range("C189").Value = 2026
$ws.Range("D189").Value = "MDPI AG"
$ws.Range("E189").Value = "Energies"
$ws.Range("F189").Value = "Kim, Woo-Jung; Lee, Yu-Seok; Chun, Yeong-Han"
$ws.Range("G189").Value = ""
$ws.Range("H189").Value = "10.3390/en19040946"
$ws.Range("I189").Value = "https://doi.org/10.3390/en19040946"
$ws.Range("J189").Value = "Journal"
$ws.Range("K189").Value = "Inverter"
$ws.Range("L189").Value = "Experiment"
$ws.Range("M189").Value = "Contacts"
$ws.Range("N189").Value = ""
$ws.Range("O189").Value = ""
$ws.Range("P189").Value = ""
$ws.Range("Q189").Value = "Robust μ-Synthesis Grid-Side Control for Inverter-Based Resources in Weak Grids"
$ws.Range("R189").Value = "High"
$ws.Range("S189").Value = "'2026-02-12"
$ws.Range("T189").Value = ""
$ws.Range("A189").Style = "Normal"
$ws.Range("G189").Style = "Normal"
$ws.Range("N189").Style = "Normal"
$ws.Range("O189").Style = "Normal"
$ws.Range("P189").Style = "Normal"
$ws.Range("S189").Style = "Normal"
$ws.Range("T189").Style = "Normal"

# Row 190
$ws.Range("A190").Value = ""
$ws.Range("B190").Value = "Experimental Investigation of Upstream Water-Level Dynamics for a Standard Open-Channel Sluice Gate and a Simplified Model"
$ws.Range("C190").Value = 2026
$ws.Range("D190").Value = "MDPI AG"
$ws.Range("E190").Value = "Water"
$ws.Range("F190").Value = "Li, Dongyan; Lv, Mouchao; Li, Hao; Jiang, Mingliang; Zhang, Wenzheng; Wang, Yingying; Qin, Jingtao"
$ws.Range("G190").Value = ""
$ws.Range("H190").Value = "10.3390/w18040476"
$ws.Range("I190").Value = "https://doi.org/10.3390/w18040476"
$ws.Range("J190").Value = "Journal"
$ws.Range("K190").Value = "n-FET"
$ws.Range("L190").Value = "Experiment"
$ws.Range("M190").Value = "Contacts"
$ws.Range("N190").Value = ""
$ws.Range("O190").Value = ""
$ws.Range("P190").Value = ""
$ws.Range("Q190").Value = "Experimental Investigation of Upstream Water-Level Dynamics for a Standard Open-Channel Sluice Gate and a Simplified Model"
$ws.Range("R190").Value = "High"
$ws.Range("S190").Value = "'2026-02-12"
$ws.Range("T190").Value = ""
$ws.Range("A190").Style = "Normal"
$ws.Range("G190").Style = "Normal"
$ws.Range("N190").Style = "Normal"
$ws.Range("O190").Style = "Normal"
$ws.Range("P190").Style = "Normal"
$ws.Range("S190").Style = "Normal"
$ws.Range("T190").Style = "Normal"

